# Add a speaker note to slide 7 ("Online format") reminding the
# presenter to check that everyone can access the chat, since this has
# caused issues before.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Accessing NotesPage lazily creates/synthesizes the notes page (image
# placeholder, notes body placeholder, slide-number placeholder) for a
# slide that doesn't have one yet.
$notesPage = $s.NotesPage

# Use the Placeholders collection (rather than a raw Shapes index) so the
# write is correctly routed to the notes body placeholder and persisted
# as this slide's notes slide part.
$notesBody = $notesPage.Shapes.Placeholders.Item(2)
$notesBody.TextFrame.TextRange.Text = "Might be worth checking if everyone can access the chat. We’ve had issues with this before."
